# act_domain_x_outcome.xlsx — "deceived v. non-deceived stats"
#
# The single existing column (header "act", 44 numeric rows) is pushed one
# column to the right (A -> B) and a brand-new leading ID column is added
# in its place: header "sub" in A1, and an integer subject id in A2:A45.
# The new ID column gets its own (non-themed) font.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Shift the existing column out of the way and create the new column A
# ---------------------------------------------------------------------
$ws.Columns("A:A").Insert()

# New header for the inserted column
$ws.Range("A1").Value2 = "sub"

# ---------------------------------------------------------------------
# 2) Fill A2:A45 with the subject ids (one 2-D array write, fast + atomic)
# ---------------------------------------------------------------------
$ids = @(
    1001,1003,1004,1006,1009,1010,1012,1013,1015,1016,1019,1021,1242,1243,1244,
    1245,1247,1248,1249,1251,1255,1276,1286,1294,1301,1302,1303,3116,3122,3125,
    3140,3143,3166,3167,3170,3173,3175,3176,3189,3190,3200,3206,3212,3220
)

$arr = New-Object 'object[,]' $ids.Length,1
for ($i = 0; $i -lt $ids.Length; $i++) {
    $arr[$i,0] = $ids[$i]
}
$ws.Range("A2:A45").Value2 = $arr

# ---------------------------------------------------------------------
# 3) Give the whole new column (header + ids) its own explicit font:
#    size 12 / black / Calibri, i.e. no longer the themed "Aptos Narrow"
# ---------------------------------------------------------------------
$idRange = $ws.Range("A1:A45")
$idRange.Font.Color = 0
$idRange.Font.Name = "Calibri"

# ---------------------------------------------------------------------
# 4) Restore the view: scrolled down to row 28, selection on E16
# ---------------------------------------------------------------------
$ws.Range("E16").Select() | Out-Null
$window = $excel.ActiveWindow
$window.ScrollRow = 28

Write-Output "done"
